$d = $word.ActiveDocument

# Locate the paragraph range to replace by content, rather than a hard-coded
# index, so the edit is resilient to unrelated shifts earlier in the doc.
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -like "*The problem here is developing a system of predicting an end count*") {
        $startIdx = $i
    }
    if ($t -like "*and the first finger (number10)*") {
        $endIdx = $i
    }
}

$start = $d.Paragraphs.Item($startIdx).Range.Start
$end = $d.Paragraphs.Item($endIdx).Range.End
$full = $d.Range($start, $end)

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = ""
$xml += "<w:p $ns><w:r><w:tab/><w:t>The problem here is developing a system of predicting an end count based on sequential system of the following:</w:t></w:r></w:p>"
$xml += "<w:p $ns/>"
$xml += "<w:p $ns><w:r><w:tab/><w:t xml:space='preserve'>• For the first five numbers, she always starts the numbering system with the thumb </w:t></w:r><w:r><w:tab/><w:t>(number 1).</w:t></w:r></w:p>"
$xml += "<w:p $ns><w:r><w:tab/><w:t xml:space='preserve'>• Followed by the first finger (number2), middle finger (number 3), ring finger (number </w:t></w:r><w:r><w:tab/><w:t>4) ending with the little finger (number 5).</w:t></w:r></w:p>"
$xml += "<w:p $ns><w:r><w:tab/><w:t xml:space='preserve'>Then the next round of 5 numbers goes like this: </w:t></w:r></w:p>"
$xml += "<w:p $ns><w:r><w:tab/><w:t xml:space='preserve'>• Ringer finger (number 6), middle finger (number 7), first finger (number 8), the thumb </w:t></w:r><w:r><w:tab/><w:t>(number 9) and the first finger (number10).</w:t></w:r></w:p>"
$xml += "<w:p $ns/>"
$xml += "<w:p $ns><w:r><w:tab/></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/><w:r><w:t>We have three main goals here: based on this sequence of instructions which finger will be identified with the ending count of 10, 100 and 1000. I believe we need to build a sequential code of the five-finger counts. Then we would implement a conditional code to change up the sequence with a set of instructions to reverse the order of count. Then back to the sequential code for a count of five. This program will last until the first total equals 10. Then we change the next program to hit a total 100. While the third program would hit a total of 1000.</w:t></w:r></w:p>"

$null = $full.InsertXML($xml)
